# edit.ps1 - Applies the commit "anado al doc cosas":
#  1) Adds a new explanatory paragraph (and a following blank paragraph) about
#     Carlos/Lidia working simultaneously, right after the paragraph that ends
#     "...Con AssertEquals de nuevo comparo el ResultadoObtenido con el
#     ResultadoEsperado." and before the "2. SUMA" heading.
#  2) Splits the run containing "...preguntar a través de una sentencia if..."
#     into two runs at "trav|és" (no visible text change).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the new paragraph + trailing blank paragraph
# ---------------------------------------------------------------------------

$anchorText = "Con AssertEquals de nuevo comparo el ResultadoObtenido con el ResultadoEsperado."

$findRng = $d.Content
$null = $findRng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The paragraph that contains the matched text is the "Con AssertEquals..." paragraph.
$anchorPara = $findRng.Paragraphs(1)
# The following (already existing, blank) paragraph.
$blankPara = $anchorPara.Next()

$insertionPoint = $blankPara.Range
$insertionPoint.Collapse(0)   # wdCollapseEnd: end of the existing blank paragraph

$newNoteText = "*Al trabajar Carlos y yo de manera simultanea nos hemos dado cuenta antes de subir mi parte a github que el también había creado la carpeta para trabajar los test y habia movido las clases a un nuevo paquete que no era el de por defecto y por tanto el proyecto me salía con conflictos. Nos hemos puesto de acuerdo y he borrado y trabajado a partir de sus modificaciones. Una vez solucionado lo he subido al repositorio local. "

# First paragraph break: starts the new paragraph that will hold $newNoteText.
$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)
$newNotePara = $insertionPoint.Paragraphs(1)

# Insert the note text into the freshly created paragraph.
$insertionPoint.InsertAfter($newNoteText)
$insertionPoint.Collapse(0)

# Second paragraph break: leaves a new blank paragraph after the note,
# right before the "2. SUMA" heading.
$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)

# The freshly-created trailing blank paragraph sometimes keeps a leftover
# formatting-only empty run after InsertParagraphAfter(). Typing a
# placeholder character and deleting it collapses that away, leaving a
# truly run-less empty paragraph (matching how the original document
# stores blank paragraphs).
$trailingBlankPara = $newNotePara.Next()
$tbRng = $trailingBlankPara.Range
$tbRng.Collapse(1)   # wdCollapseStart
$tbRng.InsertAfter("X")
$cleanupRng = $d.Range($tbRng.Start, $tbRng.Start + 1)
$cleanupRng.Delete()

# ---------------------------------------------------------------------------
# 2) Split the run "...a través de una sentencia if..." into "a trav" + "és de..."
# ---------------------------------------------------------------------------

$splitFindRng = $d.Content
$null = $splitFindRng.Find.Execute("a trav", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPoint = $splitFindRng.End
$containingPara = $d.Range($splitPoint, $splitPoint).Paragraphs(1)
$afterSplitRange = $d.Range($splitPoint, $containingPara.Range.End)

# Toggling a character property and reverting it forces Word to materialize a
# run boundary at $splitPoint without altering the visible formatting/text.
$afterSplitRange.Font.Bold = 1
$afterSplitRange.Font.Bold = 0
